$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 337 (shifts existing rows 337:443 down to 338:444)
$ws.Rows.Item(337).Insert()

# Populate the newly inserted row 337 with the new daily record
$ws.Range("A337").Value = 3
$ws.Range("B337").Value = "Femacal de La Calera"
$ws.Range("C337").Value = "Coquimbo"
$ws.Range("D337").Value = 44876
$ws.Range("E337").Value = 5
$ws.Range("F337").Value = 100112043
$ws.Range("G337").Value = "Pepino ensalada"
$ws.Range("H337").Value = "Sin especificar"
$ws.Range("I337").Value = "Primera"
$ws.Range("J337").Value = 93
$ws.Range("K337").Value = 17000
$ws.Range("L337").Value = 18000
$ws.Range("M337").Value = 17516
$ws.Range("N337").Value = "$/caja 70 unidades"
$ws.Range("O337").Value = "Región de Arica y Parinacota"
$ws.Range("P337").Value = 250
$ws.Range("Q337").Value = 70
$ws.Range("R337").Value = "Hortaliza"
